# Rewrite the monthly index data (rows 2-68) to match the refreshed source pull.
# The source re-export re-orders months within each year (Oct/Nov/Dec of the
# *previous* pull now precede Jan-Sep) and appends newly published months
# (2022 full year + 2023-01..2023-07).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure new rows below the old A1:E49 used-range inherit the same
# border/bold/centered style as the existing date column (s="1" in the OOXML),
# by priming them with a format copy from A2 before we overwrite values.
$ws.Range("A2").Copy()
$ws.Range("A50:A68").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$arr = New-Object 'object[,]' 67,5

$arr[0,0] = "2018-10"
$arr[0,1] = [double]100
$arr[0,2] = [double]97.90000000000001
$arr[0,3] = [double]100
$arr[0,4] = [double]102.1
$arr[1,0] = "2018-11"
$arr[1,1] = [double]100
$arr[1,2] = [double]98.09999999999999
$arr[1,3] = [double]100
$arr[1,4] = [double]102
$arr[2,0] = "2018-12"
$arr[2,1] = [double]100
$arr[2,2] = [double]98.59999999999999
$arr[2,3] = [double]100
$arr[2,4] = [double]102.1
$arr[3,0] = "2018-01"
$arr[3,1] = [double]100
$arr[3,2] = [double]101.4073
$arr[3,3] = [double]100
$arr[3,4] = [double]101.9503
$arr[4,0] = "2018-02"
$arr[4,1] = [double]100
$arr[4,2] = [double]100.5
$arr[4,3] = [double]100
$arr[4,4] = [double]102.5
$arr[5,0] = "2018-03"
$arr[5,1] = [double]100
$arr[5,2] = [double]99.09999999999999
$arr[5,3] = [double]100
$arr[5,4] = [double]102.4
$arr[6,0] = "2018-04"
$arr[6,1] = [double]100
$arr[6,2] = [double]98.59999999999999
$arr[6,3] = [double]100
$arr[6,4] = [double]102.4
$arr[7,0] = "2018-05"
$arr[7,1] = [double]100
$arr[7,2] = [double]98.3
$arr[7,3] = [double]100
$arr[7,4] = [double]102.5
$arr[8,0] = "2018-06"
$arr[8,1] = [double]100
$arr[8,2] = [double]97.90000000000001
$arr[8,3] = [double]100
$arr[8,4] = [double]102.4
$arr[9,0] = "2018-07"
$arr[9,1] = [double]100
$arr[9,2] = [double]97.40000000000001
$arr[9,3] = [double]100
$arr[9,4] = [double]102.3
$arr[10,0] = "2018-08"
$arr[10,1] = [double]100
$arr[10,2] = [double]98.09999999999999
$arr[10,3] = [double]100
$arr[10,4] = [double]102.5
$arr[11,0] = "2018-09"
$arr[11,1] = [double]100
$arr[11,2] = [double]97.8
$arr[11,3] = [double]100
$arr[11,4] = [double]102.2
$arr[12,0] = "2019-10"
$arr[12,1] = [double]100
$arr[12,2] = [double]105.8
$arr[12,3] = [double]100
$arr[12,4] = [double]101
$arr[13,0] = "2019-11"
$arr[13,1] = [double]100
$arr[13,2] = [double]106
$arr[13,3] = [double]100
$arr[13,4] = [double]101.1
$arr[14,0] = "2019-12"
$arr[14,1] = [double]100
$arr[14,2] = [double]105.7
$arr[14,3] = [double]100
$arr[14,4] = [double]101
$arr[15,0] = "2019-01"
$arr[15,1] = [double]100
$arr[15,2] = [double]102.2
$arr[15,3] = [double]100
$arr[15,4] = [double]102.1
$arr[16,0] = "2019-02"
$arr[16,1] = [double]100
$arr[16,2] = [double]103.3
$arr[16,3] = [double]100
$arr[16,4] = [double]101.4
$arr[17,0] = "2019-03"
$arr[17,1] = [double]100
$arr[17,2] = [double]104
$arr[17,3] = [double]100
$arr[17,4] = [double]101.7
$arr[18,0] = "2019-04"
$arr[18,1] = [double]100
$arr[18,2] = [double]104.6
$arr[18,3] = [double]100
$arr[18,4] = [double]101.7
$arr[19,0] = "2019-05"
$arr[19,1] = [double]100
$arr[19,2] = [double]104.7
$arr[19,3] = [double]100
$arr[19,4] = [double]101.6
$arr[20,0] = "2019-06"
$arr[20,1] = [double]100
$arr[20,2] = [double]105.3
$arr[20,3] = [double]100
$arr[20,4] = [double]101.5
$arr[21,0] = "2019-07"
$arr[21,1] = [double]100
$arr[21,2] = [double]106.2
$arr[21,3] = [double]100
$arr[21,4] = [double]101.5
$arr[22,0] = "2019-08"
$arr[22,1] = [double]100
$arr[22,2] = [double]105.3
$arr[22,3] = [double]100
$arr[22,4] = [double]101.1
$arr[23,0] = "2019-09"
$arr[23,1] = [double]100
$arr[23,2] = [double]105.5
$arr[23,3] = [double]100
$arr[23,4] = [double]101.1
$arr[24,0] = "2020-10"
$arr[24,1] = [double]100
$arr[24,2] = [double]101.2
$arr[24,3] = [double]100
$arr[24,4] = [double]99.59999999999999
$arr[25,0] = "2020-11"
$arr[25,1] = [double]100
$arr[25,2] = [double]101
$arr[25,3] = [double]100
$arr[25,4] = [double]99.5
$arr[26,0] = "2020-12"
$arr[26,1] = [double]100
$arr[26,2] = [double]101.3
$arr[26,3] = [double]100
$arr[26,4] = [double]99.5
$arr[27,0] = "2020-01"
$arr[27,1] = [double]100
$arr[27,2] = [double]102.6
$arr[27,3] = [double]100
$arr[27,4] = [double]100.8
$arr[28,0] = "2020-02"
$arr[28,1] = [double]100
$arr[28,2] = [double]102
$arr[28,3] = [double]100
$arr[28,4] = [double]100.6
$arr[29,0] = "2020-03"
$arr[29,1] = [double]100
$arr[29,2] = [double]101.6
$arr[29,3] = [double]100
$arr[29,4] = [double]99.8
$arr[30,0] = "2020-04"
$arr[30,1] = [double]100
$arr[30,2] = [double]101.3
$arr[30,3] = [double]100
$arr[30,4] = [double]99.5
$arr[31,0] = "2020-05"
$arr[31,1] = [double]100
$arr[31,2] = [double]101.5
$arr[31,3] = [double]100
$arr[31,4] = [double]99.40000000000001
$arr[32,0] = "2020-06"
$arr[32,1] = [double]100
$arr[32,2] = [double]101.4
$arr[32,3] = [double]100
$arr[32,4] = [double]99.40000000000001
$arr[33,0] = "2020-07"
$arr[33,1] = [double]100
$arr[33,2] = [double]101.1
$arr[33,3] = [double]100
$arr[33,4] = [double]99.40000000000001
$arr[34,0] = "2020-08"
$arr[34,1] = [double]100
$arr[34,2] = [double]101
$arr[34,3] = [double]100
$arr[34,4] = [double]99.5
$arr[35,0] = "2020-09"
$arr[35,1] = [double]100
$arr[35,2] = [double]101
$arr[35,3] = [double]100
$arr[35,4] = [double]99.5
$arr[36,0] = "2021-10"
$arr[36,1] = [double]99.90000000000001
$arr[36,2] = [double]101.9
$arr[36,3] = [double]100
$arr[36,4] = [double]101.3
$arr[37,0] = "2021-11"
$arr[37,1] = [double]99.90000000000001
$arr[37,2] = [double]101.9
$arr[37,3] = [double]100
$arr[37,4] = [double]101.1
$arr[38,0] = "2021-12"
$arr[38,1] = [double]99.90000000000001
$arr[38,2] = [double]102.4
$arr[38,3] = [double]100
$arr[38,4] = [double]101.1
$arr[39,0] = "2021-01"
$arr[39,1] = [double]100
$arr[39,2] = [double]100.8
$arr[39,3] = [double]100
$arr[39,4] = [double]99.90000000000001
$arr[40,0] = "2021-02"
$arr[40,1] = [double]99.90000000000001
$arr[40,2] = [double]100.7
$arr[40,3] = [double]100
$arr[40,4] = [double]100.4
$arr[41,0] = "2021-03"
$arr[41,1] = [double]99.90000000000001
$arr[41,2] = [double]100.6
$arr[41,3] = [double]100
$arr[41,4] = [double]100.9
$arr[42,0] = "2021-04"
$arr[42,1] = [double]99.90000000000001
$arr[42,2] = [double]100.6
$arr[42,3] = [double]100
$arr[42,4] = [double]101.2
$arr[43,0] = "2021-05"
$arr[43,1] = [double]99.90000000000001
$arr[43,2] = [double]101
$arr[43,3] = [double]100
$arr[43,4] = [double]101.2
$arr[44,0] = "2021-06"
$arr[44,1] = [double]99.90000000000001
$arr[44,2] = [double]102
$arr[44,3] = [double]100
$arr[44,4] = [double]101.2
$arr[45,0] = "2021-07"
$arr[45,1] = [double]99.90000000000001
$arr[45,2] = [double]102.1
$arr[45,3] = [double]100
$arr[45,4] = [double]101.2
$arr[46,0] = "2021-08"
$arr[46,1] = [double]99.90000000000001
$arr[46,2] = [double]102
$arr[46,3] = [double]100
$arr[46,4] = [double]101.1
$arr[47,0] = "2021-09"
$arr[47,1] = [double]99.90000000000001
$arr[47,2] = [double]102.1
$arr[47,3] = [double]100
$arr[47,4] = [double]101
$arr[48,0] = "2022-10"
$arr[48,1] = [double]99.90000000000001
$arr[48,2] = [double]102.1
$arr[48,3] = [double]100
$arr[48,4] = [double]100.6
$arr[49,0] = "2022-11"
$arr[49,1] = [double]99.90000000000001
$arr[49,2] = [double]101.5
$arr[49,3] = [double]100
$arr[49,4] = [double]100.8
$arr[50,0] = "2022-12"
$arr[50,1] = [double]99.90000000000001
$arr[50,2] = [double]100.9
$arr[50,3] = [double]100
$arr[50,4] = [double]100.8
$arr[51,0] = "2022-01"
$arr[51,1] = [double]99.90000000000001
$arr[51,2] = [double]102.5
$arr[51,3] = [double]100
$arr[51,4] = [double]101.6
$arr[52,0] = "2022-02"
$arr[52,1] = [double]99.90000000000001
$arr[52,2] = [double]102.4
$arr[52,3] = [double]100
$arr[52,4] = [double]101.2
$arr[53,0] = "2022-03"
$arr[53,1] = [double]99.90000000000001
$arr[53,2] = [double]102.8
$arr[53,3] = [double]100
$arr[53,4] = [double]101
$arr[54,0] = "2022-04"
$arr[54,1] = [double]99.90000000000001
$arr[54,2] = [double]103.2
$arr[54,3] = [double]100
$arr[54,4] = [double]100.9
$arr[55,0] = "2022-05"
$arr[55,1] = [double]99.90000000000001
$arr[55,2] = [double]102.7
$arr[55,3] = [double]100
$arr[55,4] = [double]100.9
$arr[56,0] = "2022-06"
$arr[56,1] = [double]99.90000000000001
$arr[56,2] = [double]101.6
$arr[56,3] = [double]100
$arr[56,4] = [double]100.9
$arr[57,0] = "2022-07"
$arr[57,1] = [double]99.90000000000001
$arr[57,2] = [double]101.6
$arr[57,3] = [double]100
$arr[57,4] = [double]101
$arr[58,0] = "2022-08"
$arr[58,1] = [double]99.90000000000001
$arr[58,2] = [double]101.6
$arr[58,3] = [double]100
$arr[58,4] = [double]101.1
$arr[59,0] = "2022-09"
$arr[59,1] = [double]99.90000000000001
$arr[59,2] = [double]101.5
$arr[59,3] = [double]100
$arr[59,4] = [double]101
$arr[60,0] = "2023-01"
$arr[60,1] = [double]99.90000000000001
$arr[60,2] = [double]100.6
$arr[60,3] = [double]100
$arr[60,4] = [double]100.2
$arr[61,0] = "2023-02"
$arr[61,1] = [double]101
$arr[61,2] = [double]100.7
$arr[61,3] = [double]100
$arr[61,4] = [double]100.2
$arr[62,0] = "2023-03"
$arr[62,1] = [double]101
$arr[62,2] = [double]100.5
$arr[62,3] = [double]100
$arr[62,4] = [double]100.4
$arr[63,0] = "2023-04"
$arr[63,1] = [double]101
$arr[63,2] = [double]100.1
$arr[63,3] = [double]100
$arr[63,4] = [double]100.4
$arr[64,0] = "2023-05"
$arr[64,1] = [double]101
$arr[64,2] = [double]100.2
$arr[64,3] = [double]100
$arr[64,4] = [double]100.5
$arr[65,0] = "2023-06"
$arr[65,1] = [double]101
$arr[65,2] = [double]100.3
$arr[65,3] = [double]100
$arr[65,4] = [double]100.5
$arr[66,0] = "2023-07"
$arr[66,1] = [double]101
$arr[66,2] = [double]100.4
$arr[66,3] = [double]100
$arr[66,4] = [double]100.5

$ws.Range("A2:E68").Value = $arr
